$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style (bold, border, centered) from an existing header cell (H1)
# to the new header cells so they match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-13
$values = @{
    2  = @(7, 8)
    3  = @(4, 5)
    4  = @(9, 9)
    5  = @(7, 8)
    6  = @(7, 9)
    7  = @(8, 8)
    8  = @(7, 8)
    9  = @(7, 7)
    10 = @(6, 8)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
